{"js": "// Replace the date line and each \"NN\u00d7NN=\" problem in the multiplication\n// table with its new value. Every source string in this document is\n// unique, so a plain body.search()/insertText(\"Replace\") pass is safe.\nconst replacements = [\n  [\"2025-03-06 Thursday\", \"2025-03-07 Friday\"],\n  [\"17\u00d765=\", \"83\u00d761=\"],\n  [\"12\u00d799=\", \"61\u00d750=\"],\n  [\"66\u00d781=\", \"29\u00d752=\"],\n  [\"71\u00d711=\", \"97\u00d713=\"],\n  [\"43\u00d715=\", \"61\u00d758=\"],\n  [\"57\u00d782=\", \"17\u00d750=\"],\n  [\"35\u00d758=\", \"59\u00d775=\"],\n  [\"39\u00d772=\", \"26\u00d790=\"],\n  [\"38\u00d765=\", \"42\u00d738=\"],\n  [\"51\u00d786=\", \"53\u00d714=\"],\n  [\"44\u00d720=\", \"52\u00d766=\"],\n  [\"33\u00d729=\", \"25\u00d748=\"],\n  [\"47\u00d717=\", \"26\u00d734=\"],\n  [\"28\u00d754=\", \"98\u00d792=\"],\n  [\"86\u00d758=\", \"45\u00d785=\"],\n  [\"25\u00d760=\", \"51\u00d744=\"],\n  [\"11\u00d797=\", \"75\u00d761=\"],\n  [\"84\u00d777=\", \"34\u00d711=\"],\n  [\"82\u00d780=\", \"62\u00d736=\"],\n  [\"29\u00d770=\", \"85\u00d790=\"],\n  [\"21\u00d778=\", \"60\u00d726=\"],\n  [\"49\u00d772=\", \"51\u00d790=\"],\n  [\"70\u00d753=\", \"12\u00d753=\"],\n  [\"50\u00d768=\", \"28\u00d798=\"],\n  [\"55\u00d786=\", \"16\u00d792=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NN\u00d7NN=\" problem in the multiplication\n# table with its new value. Every source string in this document is\n# unique, so a straightforward Find/Replace (wdReplaceAll) pass per pair\n# is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-06 Thursday\", \"2025-03-07 Friday\"),\n    @(\"17\u00d765=\", \"83\u00d761=\"),\n    @(\"12\u00d799=\", \"61\u00d750=\"),\n    @(\"66\u00d781=\", \"29\u00d752=\"),\n    @(\"71\u00d711=\", \"97\u00d713=\"),\n    @(\"43\u00d715=\", \"61\u00d758=\"),\n    @(\"57\u00d782=\", \"17\u00d750=\"),\n    @(\"35\u00d758=\", \"59\u00d775=\"),\n    @(\"39\u00d772=\", \"26\u00d790=\"),\n    @(\"38\u00d765=\", \"42\u00d738=\"),\n    @(\"51\u00d786=\", \"53\u00d714=\"),\n    @(\"44\u00d720=\", \"52\u00d766=\"),\n    @(\"33\u00d729=\", \"25\u00d748=\"),\n    @(\"47\u00d717=\", \"26\u00d734=\"),\n    @(\"28\u00d754=\", \"98\u00d792=\"),\n    @(\"86\u00d758=\", \"45\u00d785=\"),\n    @(\"25\u00d760=\", \"51\u00d744=\"),\n    @(\"11\u00d797=\", \"75\u00d761=\"),\n    @(\"84\u00d777=\", \"34\u00d711=\"),\n    @(\"82\u00d780=\", \"62\u00d736=\"),\n    @(\"29\u00d770=\", \"85\u00d790=\"),\n    @(\"21\u00d778=\", \"60\u00d726=\"),\n    @(\"49\u00d772=\", \"51\u00d790=\"),\n    @(\"70\u00d753=\", \"12\u00d753=\"),\n    @(\"50\u00d768=\", \"28\u00d798=\"),\n    @(\"55\u00d786=\", \"16\u00d792=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
